$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59 - Kazajistan: update daily stats
$ws.Range("B59").Value = 5126
$ws.Range("C59").Value = 36
$ws.Range("E59").Value = 3154

# Row 69 - Tailandia: update daily stats
$ws.Range("B69").Value = 3015
$ws.Range("C69").Value = 6
$ws.Range("D69").Value = 2796
$ws.Range("E69").Value = 163

# Row 98 - Kirguistan: update daily stats
$ws.Range("B98").Value = 1016
$ws.Range("C98").Value = 14
$ws.Range("D98").Value = 688
$ws.Range("E98").Value = 316

# Rows 100-102: "El Salvador" now sorts ahead of "Letonia", pushing
# Letonia and Republica de Chipre down a row, each with refreshed stats
$ws.Range("A100").Value = "El Salvador"
$ws.Range("B100").Value = 958
$ws.Range("C100").Value = 69
$ws.Range("D100").Value = 325
$ws.Range("E100").Value = 616
$ws.Range("F100").Value = 4
$ws.Range("H100").Value = 17

$ws.Range("A101").Value = "Letonia"
$ws.Range("B101").Value = 939
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 464
$ws.Range("E101").Value = 457
$ws.Range("F101").Value = 2
$ws.Range("H101").Value = 18

$ws.Range("A102").Value = "Republica de Chipre"
$ws.Range("B102").Value = 898
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 401
$ws.Range("E102").Value = 481
$ws.Range("F102").Value = 10
$ws.Range("H102").Value = 16
